# Generate Report for Handoff
#
# The 72b21453-9460-4dd3-b944-2e553c742a9f item got a fresh handoff pass,
# so its "Latest Handoff" timestamps are bumped on the Overview sheet and
# on each per-locale sheet (row 6, which is the row for that file).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the
# 72b21453-9460-4dd3-b944-2e553c742a9f.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-32-12 18:32:02"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the same file's row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-12 18:31:58"

# de-de sheet: "Latest Handoff Datetime" column (E) for the same file's row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-12 18:32:02"
